$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.012.90'
$ws.Range("E2").Value = '  -4.24%  '

$ws.Range("D3").Value = '3.797.08'
$ws.Range("E3").Value = '  -5.20%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.99'
$ws.Range("E5").Value = '  -0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.65'
$ws.Range("E6").Value = '  +1.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.662'
$ws.Range("E7").Value = '  -3.65%  '

$ws.Range("E8").Value = '  +0.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("E9").Value = '  -2.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.172'
$ws.Range("E10").Value = '  +1.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.20'
$ws.Range("E11").Value = '  -4.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000315'
$ws.Range("E12").Value = '  -1.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.14'
$ws.Range("E13").Value = '  +1.31%  '

$ws.Range("D14").Value = '4.419.74'
$ws.Range("E14").Value = '  -4.69%  '

$ws.Range("D15").Value = '3.832.74'
$ws.Range("E15").Value = '  -4.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.64'
$ws.Range("E16").Value = '  +0.81%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.58'
$ws.Range("E17").Value = '  -4.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.19'
$ws.Range("E18").Value = '  -6.67%  '

$ws.Range("E19").Value = '  -2.57%  '

$ws.Range("D20").Value = '70.078.10'
$ws.Range("E20").Value = '  -3.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '430.67'
$ws.Range("E21").Value = '  -1.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.67'
$ws.Range("E22").Value = '  -2.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '92.55'
$ws.Range("E23").Value = '  -4.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.22'
$ws.Range("E24").Value = '  -6.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.68'
$ws.Range("E25").Value = '  -4.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.00'
$ws.Range("E26").Value = '  -3.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.90'
$ws.Range("E27").Value = '  -13.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.94'
$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.27'
$ws.Range("E29").Value = '  -2.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.60'
$ws.Range("E30").Value = '  -5.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.97'
$ws.Range("E31").Value = '  +0.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.26'
$ws.Range("E32").Value = '  -3.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '47.54'
$ws.Range("E33").Value = '  -2.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.123'
$ws.Range("E34").Value = '  -6.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '69.10'
$ws.Range("E35").Value = '  -2.09%  '

$ws.Range("D36").Value = '0.0₃0965'
$ws.Range("E36").Value = '  +9.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '627.28'
$ws.Range("E37").Value = '  -6.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.417'
$ws.Range("E38").Value = '  -5.61%  '

$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.143'
$ws.Range("E41").Value = '  -2.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.21'
$ws.Range("E42").Value = '  -5.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.15'
$ws.Range("E43").Value = '  +19.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("E44").Value = '  -5.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.68'
$ws.Range("E45").Value = '  +1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.78'
$ws.Range("E46").Value = '  -9.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.142'
$ws.Range("E47").Value = '  -5.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("E48").Value = '  -15.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.23'
$ws.Range("E49").Value = '  -4.72%  '

$ws.Range("D50").Value = '2.803.51'
$ws.Range("E50").Value = '  -1.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000268'
$ws.Range("E51").Value = '  -0.52%  '
